# Update cryptos list (GitHub Actions scheduled data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "27.885.29"
$ws.Range("E2").Value = "  +4.90%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.774.86"
$ws.Range("E3").Value = "  +3.40%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.001"
$ws.Range("E4").Value = "  +0.24%  "

# Row 5 - BNB
Set-TextValue "D5" "243.01"
$ws.Range("E5").Value = "  +1.03%  "

# Row 6 - USDC
Set-TextValue "D6" "1.002"
$ws.Range("E6").Value = "  +0.23%  "

# Row 7 - XRP
Set-TextValue "D7" "0.4877"
$ws.Range("E7").Value = "  -0.74%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.2649"
$ws.Range("E8").Value = "  +2.02%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.06229"
$ws.Range("E9").Value = "  +0.42%  "

# Row 10 - WrappedEther
Set-TextValue "D10" "1.780.18"
$ws.Range("E10").Value = "  +3.58%  "

# Row 11 - Solana
Set-TextValue "D11" "16.28"
$ws.Range("E11").Value = "  +3.80%  "

# Row 12 - TRON
Set-TextValue "D12" "0.07005"
$ws.Range("E12").Value = "  +0.13%  "

# Row 13 - Polygon
Set-TextValue "D13" "0.6153"
$ws.Range("E13").Value = "  +1.59%  "

# Row 14 - Polkadot
Set-TextValue "D14" "4.594"
$ws.Range("E14").Value = "  +2.73%  "

# Row 15 - Litecoin
Set-TextValue "D15" "79.44"
$ws.Range("E15").Value = "  +3.54%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "27.891.38"
$ws.Range("E16").Value = "  +5.47%  "

# Row 17 - Dai
Set-TextValue "D17" "0.9986"
$ws.Range("E17").Value = "  -0.15%  "

# Row 18 - BinanceUSD
Set-TextValue "D18" "1.001"
$ws.Range("E18").Value = "  +0.21%  "

# Row 19 - ShibaInu
Set-TextValue "D19" "0.000007194"
$ws.Range("E19").Value = "  +0.63%  "

# Row 20 - Avalanche
Set-TextValue "D20" "11.78"
$ws.Range("E20").Value = "  +3.86%  "

# Row 21 - WrappedliquidstakedEther2.0 (E unchanged)
Set-TextValue "D21" "2.013.18"

# Row 22 - Uniswap (E unchanged)
Set-TextValue "D22" "4.562"

# Row 23 - Cosmos
Set-TextValue "D23" "8.630"
$ws.Range("E23").Value = "  +1.79%  "

# Row 24 - Chainlink
Set-TextValue "D24" "5.178"
$ws.Range("E24").Value = "  +1.83%  "

# Row 25 - Monero
Set-TextValue "D25" "141.77"
$ws.Range("E25").Value = "  +2.90%  "

# Row 26 - EthereumClassic (D unchanged)
$ws.Range("E26").Value = "  +1.98%  "

# Row 27 - LidoDAOToken
Set-TextValue "D27" "1.853"
$ws.Range("E27").Value = "  +6.48%  "

# Rows 28/29 swap: Toncoin now ranks above BitcoinCash
Set-TextValue "B28" "Toncoin"
Set-TextValue "C28" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D28" "1.417"
$ws.Range("E28").Value = "  -1.26%  "

Set-TextValue "B29" "BitcoinCash"
Set-TextValue "C29" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D29" "108.92"
$ws.Range("E29").Value = "  +2.86%  "

# Row 30 - InternetComputer(DFINITY)
Set-TextValue "D30" "4.071"
$ws.Range("E30").Value = "  +4.15%  "

# Row 31 - Stellar
Set-TextValue "D31" "0.08235"
$ws.Range("E31").Value = "  +3.78%  "

# Row 32 - Filecoin
Set-TextValue "D32" "3.755"
$ws.Range("E32").Value = "  +3.46%  "

# Row 33 - Hedera (D unchanged)
$ws.Range("E33").Value = "  +5.05%  "

# Row 34 - ARBITRUM
Set-TextValue "D34" "1.051"
$ws.Range("E34").Value = "  +5.33%  "

# Row 35 - HuobiToken
Set-TextValue "D35" "2.598"
$ws.Range("E35").Value = "  -1.14%  "

# Row 36 - ImmutableX
Set-TextValue "D36" "0.6350"
$ws.Range("E36").Value = "  +1.73%  "

# Row 37 - TrustWalletToken
Set-TextValue "D37" "0.9376"
$ws.Range("E37").Value = "  -0.65%  "

# Row 38 - MXToken
Set-TextValue "D38" "2.585"
$ws.Range("E38").Value = "  +7.15%  "

# Row 39 - RenderToken
Set-TextValue "D39" "2.043"
$ws.Range("E39").Value = "  +1.83%  "

# Row 40 - FraxShare
Set-TextValue "D40" "5.850"
$ws.Range("E40").Value = "  +6.14%  "

# Row 41 - VeChain
Set-TextValue "D41" "0.01532"
$ws.Range("E41").Value = "  +2.06%  "

# Row 42 - PaxDollar
Set-TextValue "D42" "1.001"
$ws.Range("E42").Value = "  +0.14%  "

# Row 43 - Quant
Set-TextValue "D43" "100.07"
$ws.Range("E43").Value = "  +0.50%  "

# Row 44 - TheSandbox
Set-TextValue "D44" "0.3926"
$ws.Range("E44").Value = "  +2.37%  "

# Row 45 - Aptos
Set-TextValue "D45" "7.139"
$ws.Range("E45").Value = "  +2.87%  "

# Row 46 - Algorand
Set-TextValue "D46" "0.1186"
$ws.Range("E46").Value = "  +2.89%  "

# Row 47 - Cronos (D unchanged)
$ws.Range("E47").Value = "  +0.58%  "

# Row 48 - EnergySwap
Set-TextValue "D48" "7.942"
$ws.Range("E48").Value = "  +1.99%  "

# Row 49 - Elrond
Set-TextValue "D49" "30.35"
$ws.Range("E49").Value = "  +0.37%  "

# Row 50 - NEARProtocol
Set-TextValue "D50" "1.273"
$ws.Range("E50").Value = "  +4.36%  "

# Row 51 - Aave
Set-TextValue "D51" "52.30"
$ws.Range("E51").Value = "  +1.75%  "
